$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "GNG_TO-16512556345458887"
$ws.Range("B2").Value = "go_stims-16512556345118978.csv"
$ws.Range("B3").Value = "GNG_stims-1651255634527889.csv"
$ws.Range("B4").Value = "go_stims-16512556345298893.csv"
$ws.Range("B5").Value = "GNG_stims-1651255634543893.csv"

# --- Sheet 2: NB ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "NB_TO-1651255636395462"
$ws.Range("B2").Value = "ZB-match_1-16512556346098952.csv"
$ws.Range("B3").Value = "OB-16512556349326031.csv"
$ws.Range("B4").Value = "ZB-match_0-16512556345578892.csv"
$ws.Range("B5").Value = "OB-16512556346988888.csv"
$ws.Range("B6").Value = "TB-16512556359919708.csv"
$ws.Range("B7").Value = "TB-16512556363744626.csv"
$ws.Range("B8").Value = "OB-1651255635500415.csv"
$ws.Range("B9").Value = "TB-16512556357199667.csv"
$ws.Range("B10").Value = "ZB-match_1-16512556346448972.csv"

# --- Sheet 3: RS ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "RS_TO-16512556364019458"

# --- Sheet 4: TOL ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "TOL_TO-16512556364688616"
$ws.Range("B2").Value = "MM_stims-16512556364268465.csv"
$ws.Range("B3").Value = "ZM_stims-16512556364052324.csv"
$ws.Range("B4").Value = "MM_stims-16512556364429445.csv"
$ws.Range("B5").Value = "ZM_stims-16512556364278476.csv"
$ws.Range("B6").Value = "MM_stims-16512556364678578.csv"
$ws.Range("B7").Value = "ZM_stims-16512556364439454.csv"

# --- Sheet 5: vSAT ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "vSAT_TO-16512556365449228"
$ws.Range("B2").Value = "vSAT_stims-16512556365146065.csv"
$ws.Range("B3").Value = "vSAT_stims-16512556365298655.csv"
$ws.Range("B4").Value = "SAT_stims-16512556364752088.csv"
$ws.Range("B5").Value = "SAT_stims-16512556364993453.csv"
